$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 36 (entirely empty separator row), shifting subsequent rows up by one.
$ws.Rows.Item(36).Delete()
